$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newly scheduled Metaculus median row:
# snapshot_date = 2025-12-28 (serial 46019), median_event_date = 2033-09-26 (serial 48848)
$row = 7

$ws.Range("A$row").Value = 46019
$ws.Range("B$row").Value = 48848

# Match the date formatting used by the existing data rows (yyyy-mm-dd)
$ws.Range("A$row`:B$row").NumberFormat = "yyyy-mm-dd"
